# Update the visit-count figures as discussed in canvas.
# Each edit below replaces only the single digit that actually changed;
# everything else the upstream diff shows is incidental run splitting /
# merging around unchanged text.
$d = $word.ActiveDocument

function Set-TrailingDigit($paragraphIndex, $oldDigit, $newDigit) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $r = $p.Range
    $text = $r.Text
    # Range.Text includes the trailing paragraph mark, so search within the
    # text minus that mark to find the last occurrence of the old digit.
    $body = $text.TrimEnd([char]13)
    $pos = $body.LastIndexOf($oldDigit)
    $start = $r.Start + $pos
    $target = $d.Range($start, $start + 1)
    $target.Text = $newDigit
}

# "...Number of times employee 41 entered the organization: 2" -> "...: 3"
Set-TrailingDigit 19 "2" "3"

# "...Employee 41 is the most frequent visitor with 2 visits." -> "... 3 visits."
Set-TrailingDigit 24 "2" "3"

# "41, 2" -> "41, 3"  (range 23 125 block)
Set-TrailingDigit 30 "2" "3"

# "121, 1" -> "121, 2"  (range 23 125 block)
Set-TrailingDigit 31 "1" "2"

# "41, 2" -> "41, 3"  (range 30 100 block)
Set-TrailingDigit 34 "2" "3"
